$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 725, shifting existing rows 725:822 down to 726:823
$ws.Rows(725).Insert()

# Populate the newly inserted row with the new record
$ws.Range("A725").Value = 3
$ws.Range("B725").Value = "Femacal de La Calera"
$ws.Range("C725").Value = "Coquimbo"
$ws.Range("D725").Value = 45154
$ws.Range("E725").Value = 5
$ws.Range("F725").Value = 100112032
$ws.Range("G725").Value = "Zapallo italiano"
$ws.Range("H725").Value = "Sin especificar"
$ws.Range("I725").Value = "Primera"
$ws.Range("J725").Value = 100
$ws.Range("K725").Value = 16000
$ws.Range("L725").Value = 17000
$ws.Range("M725").Value = 16500
$ws.Range("N725").Value = "$/caja 60 unidades"
$ws.Range("O725").Value = "Región de Arica y Parinacota"
$ws.Range("P725").Value = 275
$ws.Range("Q725").Value = 60
$ws.Range("R725").Value = "Hortaliza"
